# Update Adam12-Sdc4 LR-pair table with new TPM-derived values and add the
# "Resolving-Mac" target-cluster rows (rows 22-26) that were missing before.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the full replacement block (rows 2-26, columns A-T) as a 2-D array
# and write it in a single shot so numeric precision round-trips exactly.
$data = New-Object 'object[,]' 25,20
# row 2
$data[0,0] = "ECs"
$data[0,1] = "Adam12"
$data[0,2] = "Sdc4"
$data[0,3] = "ECs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 4.480363666666666
$data[0,7] = 13.441091
$data[0,8] = 0.05823429740900917
$data[0,9] = 0.05886574272937452
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 1.378421333333333
$data[0,13] = 4.135264
$data[0,14] = 0.01656231489052403
$data[0,15] = 0.01794267551419991
$data[0,16] = 6.17582885922489
$data[0,17] = 55.582459733024
$data[0,18] = 0.0009644947711164376
$data[0,19] = 0.00105620892069554

# row 3
$data[1,0] = "ECs"
$data[1,1] = "Adam12"
$data[1,2] = "Sdc4"
$data[1,3] = "FAPs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 4.480363666666666
$data[1,7] = 13.441091
$data[1,8] = 0.05823429740900917
$data[1,9] = 0.05886574272937452
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 18.067884
$data[1,13] = 54.20365200000001
$data[1,14] = 0.2170932623988173
$data[1,15] = 0.2351865659654651
$data[1,16] = 80.95069100714801
$data[1,17] = 728.556219064332
$data[1,18] = 0.01264227360802479
$data[1,19] = 0.01384443188552814

# row 4
$data[2,0] = "ECs"
$data[2,1] = "Adam12"
$data[2,2] = "Sdc4"
$data[2,3] = "Inflammatory-Mac"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 4.480363666666666
$data[2,7] = 13.441091
$data[2,8] = 0.05823429740900917
$data[2,9] = 0.05886574272937452
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 17.58286933333333
$data[2,13] = 52.748608
$data[2,14] = 0.2112656061941426
$data[2,15] = 0.22887321273073
$data[2,16] = 78.77764891681421
$data[2,17] = 708.9988402513279
$data[2,18] = 0.01230290414340431
$data[2,19] = 0.01347279165825256

# row 5
$data[3,0] = "ECs"
$data[3,1] = "Adam12"
$data[3,2] = "Sdc4"
$data[3,3] = "MuSCs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 4.480363666666666
$data[3,7] = 13.441091
$data[3,8] = 0.05823429740900917
$data[3,9] = 0.05886574272937452
$data[3,10] = 2
$data[3,11] = 1
$data[3,12] = 19.2082395
$data[3,13] = 38.416479
$data[3,14] = 0.2307951156866419
$data[3,15] = 0.1666869194070983
$data[3,16] = 86.0598983564315
$data[3,17] = 516.359390138589
$data[3,18] = 0.01344019140744258
$data[3,19] = 0.009812149314170236

# row 6
$data[4,0] = "ECs"
$data[4,1] = "Adam12"
$data[4,2] = "Sdc4"
$data[4,3] = "Resolving-Mac"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 4.480363666666666
$data[4,7] = 13.441091
$data[4,8] = 0.05823429740900917
$data[4,9] = 0.05886574272937452
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 26.988955
$data[4,13] = 80.966865
$data[4,14] = 0.3242837008298742
$data[4,15] = 0.3513106263825066
$data[4,16] = 120.9203333833017
$data[4,17] = 1088.283000449715
$data[4,18] = 0.01888443347902105
$data[4,19] = 0.02068016095072805

# row 7
$data[5,0] = "FAPs"
$data[5,1] = "Adam12"
$data[5,2] = "Sdc4"
$data[5,3] = "ECs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 69.95512000000001
$data[5,7] = 209.86536
$data[5,8] = 0.90925370493279
$data[5,9] = 0.9191129120074827
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 1.378421333333333
$data[5,13] = 4.135264
$data[5,14] = 0.01656231489052403
$data[5,15] = 0.01794267551419991
$data[5,16] = 96.42762978389335
$data[5,17] = 867.8486680550401
$data[5,18] = 0.01505934617647249
$data[5,19] = 0.01649134474106164

# row 8
$data[6,0] = "FAPs"
$data[6,1] = "Adam12"
$data[6,2] = "Sdc4"
$data[6,3] = "FAPs"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 69.95512000000001
$data[6,7] = 209.86536
$data[6,8] = 0.90925370493279
$data[6,9] = 0.9191129120074827
$data[6,10] = 3
$data[6,11] = 1
$data[6,12] = 18.067884
$data[6,13] = 54.20365200000001
$data[6,14] = 0.2170932623988173
$data[6,15] = 0.2351865659654651
$data[6,16] = 1263.94099336608
$data[6,17] = 11375.46894029472
$data[6,18] = 0.1973928531520709
$data[6,19] = 0.2161630095095586

# row 9
$data[7,0] = "FAPs"
$data[7,1] = "Adam12"
$data[7,2] = "Sdc4"
$data[7,3] = "Inflammatory-Mac"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 69.95512000000001
$data[7,7] = 209.86536
$data[7,8] = 0.90925370493279
$data[7,9] = 0.9191129120074827
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 17.58286933333333
$data[7,13] = 52.748608
$data[7,14] = 0.2112656061941426
$data[7,15] = 0.22887321273073
$data[7,16] = 1230.011734157653
$data[7,17] = 11070.10560741888
$data[7,18] = 0.1920940351568959
$data[7,19] = 0.2103603250334493

# row 10
$data[8,0] = "FAPs"
$data[8,1] = "Adam12"
$data[8,2] = "Sdc4"
$data[8,3] = "MuSCs"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 69.95512000000001
$data[8,7] = 209.86536
$data[8,8] = 0.90925370493279
$data[8,9] = 0.9191129120074827
$data[8,10] = 2
$data[8,11] = 1
$data[8,12] = 19.2082395
$data[8,13] = 38.416479
$data[8,14] = 0.2307951156866419
$data[8,15] = 0.1666869194070983
$data[8,16] = 1343.71469921124
$data[8,17] = 8062.288195267441
$data[8,18] = 0.209851314018471
$data[8,19] = 0.1532040998898147

# row 11
$data[9,0] = "FAPs"
$data[9,1] = "Adam12"
$data[9,2] = "Sdc4"
$data[9,3] = "Resolving-Mac"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 69.95512000000001
$data[9,7] = 209.86536
$data[9,8] = 0.90925370493279
$data[9,9] = 0.9191129120074827
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 26.988955
$data[9,13] = 80.966865
$data[9,14] = 0.3242837008298742
$data[9,15] = 0.3513106263825066
$data[9,16] = 1888.0155856996
$data[9,17] = 16992.1402712964
$data[9,18] = 0.2948561564288796
$data[9,19] = 0.3228941328335984

# row 12
$data[10,0] = "Inflammatory-Mac"
$data[10,1] = "Adam12"
$data[10,2] = "Sdc4"
$data[10,3] = "ECs"
$data[10,4] = 1
$data[10,5] = 0.3333333333333333
$data[10,6] = 0.02020466666666667
$data[10,7] = 0.060614
$data[10,8] = 0.0002626136303332581
$data[10,9] = 0.0002654611987820265
$data[10,10] = 3
$data[10,11] = 1
$data[10,12] = 1.378421333333333
$data[10,13] = 4.135264
$data[10,14] = 0.01656231489052403
$data[10,15] = 0.01794267551419991
$data[10,16] = 0.02785054356622222
$data[10,17] = 0.250654892096
$data[10,18] = [double]"4.349489640123094e-06"
$data[10,19] = [double]"4.763084151356423e-06"

# row 13
$data[11,0] = "Inflammatory-Mac"
$data[11,1] = "Adam12"
$data[11,2] = "Sdc4"
$data[11,3] = "FAPs"
$data[11,4] = 1
$data[11,5] = 0.3333333333333333
$data[11,6] = 0.02020466666666667
$data[11,7] = 0.060614
$data[11,8] = 0.0002626136303332581
$data[11,9] = 0.0002654611987820265
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 18.067884
$data[11,13] = 54.20365200000001
$data[11,14] = 0.2170932623988173
$data[11,15] = 0.2351865659654651
$data[11,16] = 0.365055573592
$data[11,17] = 3.285500162328
$data[11,18] = [double]"5.701164975944399e-05"
$data[11,19] = [double]"6.243290773862053e-05"

# row 14
$data[12,0] = "Inflammatory-Mac"
$data[12,1] = "Adam12"
$data[12,2] = "Sdc4"
$data[12,3] = "Inflammatory-Mac"
$data[12,4] = 1
$data[12,5] = 0.3333333333333333
$data[12,6] = 0.02020466666666667
$data[12,7] = 0.060614
$data[12,8] = 0.0002626136303332581
$data[12,9] = 0.0002654611987820265
$data[12,10] = 3
$data[12,11] = 1
$data[12,12] = 17.58286933333333
$data[12,13] = 52.748608
$data[12,14] = 0.2112656061941426
$data[12,15] = 0.22887321273073
$data[12,16] = 0.3552560139235555
$data[12,17] = 3.197304125312
$data[12,18] = [double]"5.548122780720023e-05"
$data[12,19] = [double]"6.075695742059337e-05"

# row 15
$data[13,0] = "Inflammatory-Mac"
$data[13,1] = "Adam12"
$data[13,2] = "Sdc4"
$data[13,3] = "MuSCs"
$data[13,4] = 1
$data[13,5] = 0.3333333333333333
$data[13,6] = 0.02020466666666667
$data[13,7] = 0.060614
$data[13,8] = 0.0002626136303332581
$data[13,9] = 0.0002654611987820265
$data[13,10] = 2
$data[13,11] = 1
$data[13,12] = 19.2082395
$data[13,13] = 38.416479
$data[13,14] = 0.2307951156866419
$data[13,15] = 0.1666869194070983
$data[13,16] = 0.388096076351
$data[13,17] = 2.328576458106
$data[13,18] = [double]"6.06099431936533e-05"
$data[13,19] = [double]"4.424890944709137e-05"

# row 16
$data[14,0] = "Inflammatory-Mac"
$data[14,1] = "Adam12"
$data[14,2] = "Sdc4"
$data[14,3] = "Resolving-Mac"
$data[14,4] = 1
$data[14,5] = 0.3333333333333333
$data[14,6] = 0.02020466666666667
$data[14,7] = 0.060614
$data[14,8] = 0.0002626136303332581
$data[14,9] = 0.0002654611987820265
$data[14,10] = 3
$data[14,11] = 1
$data[14,12] = 26.988955
$data[14,13] = 80.966865
$data[14,14] = 0.3242837008298742
$data[14,15] = 0.3513106263825066
$data[14,16] = 0.5453028394566667
$data[14,17] = 4.90772555511
$data[14,18] = [double]"8.516131993283743e-05"
$data[14,19] = [double]"9.325934002436483e-05"

# row 17
$data[15,0] = "MuSCs"
$data[15,1] = "Adam12"
$data[15,2] = "Sdc4"
$data[15,3] = "ECs"
$data[15,4] = 2
$data[15,5] = 1
$data[15,6] = 2.4758755
$data[15,7] = 4.951751
$data[15,8] = 0.03218061767783864
$data[15,9] = 0.02168637206800571
$data[15,10] = 3
$data[15,11] = 1
$data[15,12] = 1.378421333333333
$data[15,13] = 4.135264
$data[15,14] = 0.01656231489052403
$data[15,15] = 0.01794267551419991
$data[15,16] = 3.412799607877333
$data[15,17] = 20.476797647264
$data[15,18] = 0.0005329855233519278
$data[15,19] = 0.0003891115370964351

# row 18
$data[16,0] = "MuSCs"
$data[16,1] = "Adam12"
$data[16,2] = "Sdc4"
$data[16,3] = "FAPs"
$data[16,4] = 2
$data[16,5] = 1
$data[16,6] = 2.4758755
$data[16,7] = 4.951751
$data[16,8] = 0.03218061767783864
$data[16,9] = 0.02168637206800571
$data[16,10] = 3
$data[16,11] = 1
$data[16,12] = 18.067884
$data[16,13] = 54.20365200000001
$data[16,14] = 0.2170932623988173
$data[16,15] = 0.2351865659654651
$data[16,16] = 44.73383133244201
$data[16,17] = 268.402987994652
$data[16,18] = 0.006986195277691042
$data[16,19] = 0.005100343374923646

# row 19
$data[17,0] = "MuSCs"
$data[17,1] = "Adam12"
$data[17,2] = "Sdc4"
$data[17,3] = "Inflammatory-Mac"
$data[17,4] = 2
$data[17,5] = 1
$data[17,6] = 2.4758755
$data[17,7] = 4.951751
$data[17,8] = 0.03218061767783864
$data[17,9] = 0.02168637206800571
$data[17,10] = 3
$data[17,11] = 1
$data[17,12] = 17.58286933333333
$data[17,13] = 52.748608
$data[17,14] = 0.2112656061941426
$data[17,15] = 0.22887321273073
$data[17,16] = 43.53299540210133
$data[17,17] = 261.197972412608
$data[17,18] = 0.006798657701410521
$data[17,19] = 0.004963429647678434

# row 20
$data[18,0] = "MuSCs"
$data[18,1] = "Adam12"
$data[18,2] = "Sdc4"
$data[18,3] = "MuSCs"
$data[18,4] = 2
$data[18,5] = 1
$data[18,6] = 2.4758755
$data[18,7] = 4.951751
$data[18,8] = 0.03218061767783864
$data[18,9] = 0.02168637206800571
$data[18,10] = 2
$data[18,11] = 1
$data[18,12] = 19.2082395
$data[18,13] = 38.416479
$data[18,14] = 0.2307951156866419
$data[18,15] = 0.1666869194070983
$data[18,16] = 47.55720957618225
$data[18,17] = 190.228838304729
$data[18,18] = 0.007427129379824363
$data[18,19] = 0.003614834553132017

# row 21
$data[19,0] = "MuSCs"
$data[19,1] = "Adam12"
$data[19,2] = "Sdc4"
$data[19,3] = "Resolving-Mac"
$data[19,4] = 2
$data[19,5] = 1
$data[19,6] = 2.4758755
$data[19,7] = 4.951751
$data[19,8] = 0.03218061767783864
$data[19,9] = 0.02168637206800571
$data[19,10] = 3
$data[19,11] = 1
$data[19,12] = 26.988955
$data[19,13] = 80.966865
$data[19,14] = 0.3242837008298742
$data[19,15] = 0.3513106263825066
$data[19,16] = 66.8212924551025
$data[19,17] = 400.927754730615
$data[19,18] = 0.01043564979556079
$data[19,19] = 0.007618652955175183

# row 22
$data[20,0] = "Resolving-Mac"
$data[20,1] = "Adam12"
$data[20,2] = "Sdc4"
$data[20,3] = "ECs"
$data[20,4] = 1
$data[20,5] = 0.3333333333333333
$data[20,6] = 0.005290666666666667
$data[20,7] = 0.015872
$data[20,8] = [double]"6.876635002886251e-05"
$data[20,9] = [double]"6.951199635510484e-05"
$data[20,10] = 3
$data[20,11] = 1
$data[20,12] = 1.378421333333333
$data[20,13] = 4.135264
$data[20,14] = 0.01656231489052403
$data[20,15] = 0.01794267551419991
$data[20,16] = 0.007292767800888889
$data[20,17] = 0.06563491020800001
$data[20,18] = [double]"1.138929943050017e-06"
$data[20,19] = [double]"1.247231194943893e-06"

# row 23
$data[21,0] = "Resolving-Mac"
$data[21,1] = "Adam12"
$data[21,2] = "Sdc4"
$data[21,3] = "FAPs"
$data[21,4] = 1
$data[21,5] = 0.3333333333333333
$data[21,6] = 0.005290666666666667
$data[21,7] = 0.015872
$data[21,8] = [double]"6.876635002886251e-05"
$data[21,9] = [double]"6.951199635510484e-05"
$data[21,10] = 3
$data[21,11] = 1
$data[21,12] = 18.067884
$data[21,13] = 54.20365200000001
$data[21,14] = 0.2170932623988173
$data[21,15] = 0.2351865659654651
$data[21,16] = 0.09559115161600001
$data[21,17] = 0.8603203645440001
$data[21,18] = [double]"1.492871127102476e-05"
$data[21,19] = [double]"1.634828771616104e-05"

# row 24
$data[22,0] = "Resolving-Mac"
$data[22,1] = "Adam12"
$data[22,2] = "Sdc4"
$data[22,3] = "Inflammatory-Mac"
$data[22,4] = 1
$data[22,5] = 0.3333333333333333
$data[22,6] = 0.005290666666666667
$data[22,7] = 0.015872
$data[22,8] = [double]"6.876635002886251e-05"
$data[22,9] = [double]"6.951199635510484e-05"
$data[22,10] = 3
$data[22,11] = 1
$data[22,12] = 17.58286933333333
$data[22,13] = 52.748608
$data[22,14] = 0.2112656061941426
$data[22,15] = 0.22887321273073
$data[22,16] = 0.0930251006862222
$data[22,17] = 0.837225906176
$data[22,18] = [double]"1.452796462460623e-05"
$data[22,19] = [double]"1.590943392911964e-05"

# row 25
$data[23,0] = "Resolving-Mac"
$data[23,1] = "Adam12"
$data[23,2] = "Sdc4"
$data[23,3] = "MuSCs"
$data[23,4] = 1
$data[23,5] = 0.3333333333333333
$data[23,6] = 0.005290666666666667
$data[23,7] = 0.015872
$data[23,8] = [double]"6.876635002886251e-05"
$data[23,9] = [double]"6.951199635510484e-05"
$data[23,10] = 2
$data[23,11] = 1
$data[23,12] = 19.2082395
$data[23,13] = 38.416479
$data[23,14] = 0.2307951156866419
$data[23,15] = 0.1666869194070983
$data[23,16] = 0.101624392448
$data[23,17] = 0.6097463546880001
$data[23,18] = [double]"1.587093771025944e-05"
$data[23,19] = [double]"1.158674053426988e-05"

# row 26
$data[24,0] = "Resolving-Mac"
$data[24,1] = "Adam12"
$data[24,2] = "Sdc4"
$data[24,3] = "Resolving-Mac"
$data[24,4] = 1
$data[24,5] = 0.3333333333333333
$data[24,6] = 0.005290666666666667
$data[24,7] = 0.015872
$data[24,8] = [double]"6.876635002886251e-05"
$data[24,9] = [double]"6.951199635510484e-05"
$data[24,10] = 3
$data[24,11] = 1
$data[24,12] = 26.988955
$data[24,13] = 80.966865
$data[24,14] = 0.3242837008298742
$data[24,15] = 0.3513106263825066
$data[24,16] = 0.1427895645866667
$data[24,17] = 1.28510608128
$data[24,18] = [double]"2.229980647992206e-05"
$data[24,19] = [double]"2.44203029806104e-05"

$ws.Range("A2:T26").Value = $data